# Auto-generated edit script applying Tonberry_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1961.3334
$ws.Range("I6").Value = 451
$ws.Range("J6").Value = 2542.2307
$ws.Range("K6").Value = 1353
$ws.Range("L6").Value = 7626.6921
$ws.Range("M6").Value = -1241
$ws.Range("N6").Value = -7850.6921
$ws.Range("H15").Value = 1005.2941
$ws.Range("I15").Value = 1005.2941
$ws.Range("K15").Value = 3015.8823
$ws.Range("M15").Value = -2846.8823
$ws.Range("H112").Value = 2859.8125
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2859.8125
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 8579.4375
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -10795.4375
$ws.Range("H135").Value = 454.76923
$ws.Range("I135").Value = 439.18918
$ws.Range("K135").Value = 3952.70262
$ws.Range("M135").Value = -1417.70262
$ws.Range("H137").Value = 31695.06
$ws.Range("I137").Value = 1209.375
$ws.Range("J137").Value = 112990.22
$ws.Range("K137").Value = 3628.125
$ws.Range("L137").Value = 338970.66
$ws.Range("M137").Value = -1078.125
$ws.Range("N137").Value = -344070.66
$ws.Range("H138").Value = 1981.2325
$ws.Range("I138").Value = 1789.8148
$ws.Range("J138").Value = 2304.25
$ws.Range("K138").Value = 5369.4444
$ws.Range("L138").Value = 6912.75
$ws.Range("M138").Value = -229.4444000000003
$ws.Range("N138").Value = -17192.75
$ws.Range("H141").Value = 623778.9399999999
$ws.Range("I141").Value = 683673.0600000001
$ws.Range("K141").Value = 2051019.18
$ws.Range("M141").Value = -2045839.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3444.2173
$ws.Range("I32").Value = 2851.4824
$ws.Range("J32").Value = 10641.714
$ws.Range("K32").Value = 2851.4824
$ws.Range("L32").Value = 10641.714
$ws.Range("M32").Value = -2564.4824
$ws.Range("N32").Value = -11215.714
$ws.Range("H61").Value = 3263.1035
$ws.Range("I61").Value = 1123.3043
$ws.Range("K61").Value = 1123.3043
$ws.Range("M61").Value = -911.3043
$ws.Range("H102").Value = 1499.8
$ws.Range("I102").Value = 1499.8
$ws.Range("K102").Value = 1499.8
$ws.Range("M102").Value = 122.2
$ws.Range("H122").Value = 1696.4073
$ws.Range("I122").Value = 1632.7142
$ws.Range("J122").Value = 1919.3334
$ws.Range("K122").Value = 4898.142599999999
$ws.Range("L122").Value = 5758.0002
$ws.Range("M122").Value = -2448.142599999999
$ws.Range("N122").Value = -10658.0002
$ws.Range("H132").Value = 1326.6279
$ws.Range("I132").Value = 1092.6052
$ws.Range("K132").Value = 3277.8156
$ws.Range("M132").Value = -747.8155999999999
$ws.Range("H136").Value = 3263.1035
$ws.Range("I136").Value = 1123.3043
$ws.Range("K136").Value = 3369.9129
$ws.Range("M136").Value = -819.9129000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 25000
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("N88").Value = -25812
$ws.Range("H91").Value = 25000
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("N91").Value = -27808
$ws.Range("H94").Value = 1632.2222
$ws.Range("I94").Value = 2018.5714
$ws.Range("J94").Value = 280
$ws.Range("K94").Value = 2018.5714
$ws.Range("L94").Value = 280
$ws.Range("M94").Value = -1567.5714
$ws.Range("N94").Value = -1182
$ws.Range("H134").Value = 6376.8335
$ws.Range("I134").Value = 6862.2
$ws.Range("J134").Value = 3950
$ws.Range("K134").Value = 20586.6
$ws.Range("L134").Value = 11850
$ws.Range("M134").Value = -18051.6
$ws.Range("N134").Value = -16920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 12268.167
$ws.Range("I32").Value = 4024.5
$ws.Range("J32").Value = 28755.5
$ws.Range("K32").Value = 4024.5
$ws.Range("L32").Value = 28755.5
$ws.Range("M32").Value = -3708.5
$ws.Range("N32").Value = -29387.5
$ws.Range("H45").Value = 14299.667
$ws.Range("I45").Value = 14299.667
$ws.Range("K45").Value = 14299.667
$ws.Range("M45").Value = -13706.667
$ws.Range("H132").Value = 1890.1364
$ws.Range("I132").Value = 1334.4242
$ws.Range("J132").Value = 3557.2727
$ws.Range("K132").Value = 4003.2726
$ws.Range("L132").Value = 10671.8181
$ws.Range("M132").Value = -1473.2726
$ws.Range("N132").Value = -15731.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4000
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13372
$ws.Range("H65").Value = 4000
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 36000
$ws.Range("N65").Value = -42864
$ws.Range("H70").Value = 1353
$ws.Range("I70").Value = 804
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 2412
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -2097
$ws.Range("N70").Value = -9630
$ws.Range("H73").Value = 1353
$ws.Range("I73").Value = 804
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 2412
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -1320
$ws.Range("N73").Value = -11184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 21399.166
$ws.Range("J92").Value = 21399.166
$ws.Range("L92").Value = 21399.166
$ws.Range("N92").Value = -25143.166
$ws.Range("H102").Value = 4042.9167
$ws.Range("I102").Value = 5047.5
$ws.Range("J102").Value = 2636.5
$ws.Range("K102").Value = 5047.5
$ws.Range("L102").Value = 2636.5
$ws.Range("M102").Value = -3425.5
$ws.Range("N102").Value = -5880.5
$ws.Range("H132").Value = 727653.5
$ws.Range("I132").Value = 1167005.2
$ws.Range("J132").Value = 2723.1
$ws.Range("K132").Value = 3501015.6
$ws.Range("L132").Value = 8169.299999999999
$ws.Range("M132").Value = -3498485.6
$ws.Range("N132").Value = -13229.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 7967.25
$ws.Range("I32").Value = 7333
$ws.Range("K32").Value = 7333
$ws.Range("M32").Value = -7016
$ws.Range("H100").Value = 2125.7144
$ws.Range("I100").Value = 1650
$ws.Range("J100").Value = 4980
$ws.Range("K100").Value = 1650
$ws.Range("L100").Value = 4980
$ws.Range("M100").Value = -1109
$ws.Range("N100").Value = -6062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 47500
$ws.Range("J68").Value = 47500
$ws.Range("L68").Value = 47500
$ws.Range("N68").Value = -49122
$ws.Range("H71").Value = 47500
$ws.Range("J71").Value = 47500
$ws.Range("L71").Value = 142500
$ws.Range("N71").Value = -150612
$ws.Range("H132").Value = 1203.196
$ws.Range("I132").Value = 721.87177
$ws.Range("J132").Value = 2767.5
$ws.Range("K132").Value = 2165.61531
$ws.Range("L132").Value = 8302.5
$ws.Range("M132").Value = 364.3846899999999
$ws.Range("N132").Value = -13362.5
$ws.Range("H136").Value = 12921641
$ws.Range("I136").Value = 15433996
$ws.Range("J136").Value = 960.1429000000001
$ws.Range("K136").Value = 46301988
$ws.Range("L136").Value = 2880.4287
$ws.Range("M136").Value = -46299438
$ws.Range("N136").Value = -7980.4287
